# add 3.3V regulator to part list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bezeichnung (name) for the new part
$ws.Range("B13").Value = "3.3V Spannungsregler"

# Mouser Link column: hyperlink whose display text is the URL itself,
# matching the style used by the other rows in this table.
$url = "https://www.mouser.de/ProductDetail/Texas-Instruments/UA78M33CDCY?qs=sbcp%2F4gpy09US8tH%252B2FxOw%3D%3D"
$ws.Hyperlinks.Add($ws.Range("C13"), $url, "", "", $url)
$ws.Range("C13").Style = "Link"

# Schematic Bezeichner (designator) for the new part
$ws.Range("A13").Value = "IC3"

# Mirrors the author's recorded selection after the edit
$ws.Range("B17").Select()
